{"js": "// 1. Append \" 1\" to the title paragraph (\"Team Buttercup Minutes\" -> \"Team Buttercup Minutes 1\").\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst titlePara = paragraphs.items[0];\ntitlePara.insertText(\" 1\", Word.InsertLocation.end);\nawait context.sync();\n\n// 2. Fix the \"Minuets\" typo -> \"Minutes\" in the meeting-info line.\nconst typoResults = context.document.body.search(\"Minuets\", { matchCase: true });\ntypoResults.load(\"text\");\nawait context.sync();\nif (typoResults.items.length > 0) {\n  typoResults.items[0].insertText(\"Minutes\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 3. Merge the \"Chief Tester: Patrick\" / \" Evens \" runs into a single run\n//    by replacing the combined text in place (formatting is identical, so\n//    this is a no-op visually but collapses the run split).\nconst chiefResults = context.document.body.search(\"Chief Tester: Patrick Evens \", { matchCase: true });\nchiefResults.load(\"text\");\nawait context.sync();\nif (chiefResults.items.length > 0) {\n  chiefResults.items[0].insertText(\"Chief Tester: Patrick Evens \", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 4. Move the \"_GoBack\" bookmark so it collapses at the end of the title\n//    paragraph instead of spanning the whole document body.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst paragraphsAfter = context.document.body.paragraphs;\nparagraphsAfter.load(\"text\");\nawait context.sync();\n\nconst titleParaAfter = paragraphsAfter.items[0];\nconst titleEndAfter = titleParaAfter.getRange(\"End\");\ntitleEndAfter.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Append \" 1\" to the title paragraph (\"Team Buttercup Minutes\" -> \"Team Buttercup Minutes 1\").\n$titlePara = $d.Paragraphs(1)\n$titleLastChar = $titlePara.Range.Characters($titlePara.Range.Characters.Count - 1)\n$titleLastChar.Collapse(0)   # wdCollapseEnd\n$titleLastChar.InsertAfter(\" 1\")\n\n# 2. Fix the \"Minuets\" typo -> \"Minutes\" in the meeting-info line.\n$findRng = $d.Content\n$findRng.Find.ClearFormatting()\n$findRng.Find.Text = \"Minuets\"\n$findRng.Find.Replacement.Text = \"Minutes\"\n$findRng.Find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n# 3. Merge the \"Chief Tester: Patrick\" / \" Evens \" runs into a single run by\n#    replacing the combined text in place (formatting is identical, so this\n#    is a no-op visually but collapses the run split).\n$chiefRng = $d.Content\n$chiefRng.Find.ClearFormatting()\n$chiefRng.Find.Text = \"Chief Tester: Patrick Evens \"\n$chiefRng.Find.Replacement.Text = \"Chief Tester: Patrick Evens \"\n$chiefRng.Find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n# 4. Move the \"_GoBack\" bookmark so it collapses at the end of the title\n#    paragraph instead of spanning the whole document body.\n#\n#    A bookmark can't be added directly at the exact paragraph-end boundary,\n#    so a temporary placeholder character is inserted there, the bookmark is\n#    anchored immediately before it (recording that position up front), and\n#    the placeholder is removed again \u2014 leaving the bookmark collapsed in\n#    the right spot.\n$bm = $d.Bookmarks(\"_GoBack\")\n$bm.Delete()\n\n$titlePara2 = $d.Paragraphs(1)\n$placeholderAnchor = $titlePara2.Range.Characters($titlePara2.Range.Characters.Count - 1)\n$placeholderAnchor.Collapse(0)   # wdCollapseEnd\n$insertPos = $placeholderAnchor.Start\n$placeholderAnchor.InsertAfter(\"X\")\n\n$bmRange = $d.Range($insertPos, $insertPos)\n$d.Bookmarks.Add(\"_GoBack\", $bmRange)\n\n$placeholderRange = $d.Range($insertPos, $insertPos + 1)\n$placeholderRange.Delete()\n"}
